$wb = $excel.ActiveWorkbook

$values = @(0.7070674194434616, -0.21945299999999968, 0.13247858013724567, -0.4540000000000002, 1.5829618029997903, 16.12947350163202, 2.67659686508821)
$cols = @("C", "D", "E", "F", "G", "H", "I")

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])11"
        $ws.Range($addr).Value = $values[$i]
    }
}
